$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet: update the "Date" property value
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Cells.Item(8, 2).Value2 = "2026-01-22T09:24:45+00:00"

# ---------------------------------------------------------------------------
# 2) Elements sheet: split "dateDebutFin" into "dateDebut" + "dateFin"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Row 7 currently describes "fr-lm-effet-indesirable.dateDebutFin".
# Rename it to describe only the "dateDebut" part.
$ws.Cells.Item(7, 1).Value2  = "fr-lm-effet-indesirable.dateDebut"   # ID
$ws.Cells.Item(7, 2).Value2  = "fr-lm-effet-indesirable.dateDebut"   # Path
$ws.Cells.Item(7, 32).Value2 = "fr-lm-effet-indesirable.dateDebut"   # Base Path
$ws.Cells.Item(7, 12).Value2 = "Date de début de l'effet indésirable" # Short
$ws.Cells.Item(7, 13).Value2 = "Date de début de l'effet indésirable" # Definition

# Insert a brand new row right after it to hold "dateFin"
$ws.Rows.Item(8).Insert()

# Force the numeric-looking text cells (Min/Max/Base Min/Base Max) to stay
# text, matching how the rest of the sheet stores "0"/"1" as text, not
# numbers.
$f8 = $ws.Cells.Item(8, 6)
$f8.NumberFormat = "@"
$f8.Value2 = "0"

$g8 = $ws.Cells.Item(8, 7)
$g8.NumberFormat = "@"
$g8.Value2 = "1"

$ag8 = $ws.Cells.Item(8, 33)
$ag8.NumberFormat = "@"
$ag8.Value2 = "0"

$ah8 = $ws.Cells.Item(8, 34)
$ah8.NumberFormat = "@"
$ah8.Value2 = "1"

# Regular text cells for the new "dateFin" row
$ws.Cells.Item(8, 1).Value2  = "fr-lm-effet-indesirable.dateFin"       # ID
$ws.Cells.Item(8, 2).Value2  = "fr-lm-effet-indesirable.dateFin"       # Path
$ws.Cells.Item(8, 11).Value2 = "dateTime`n"                             # Type(s)
$ws.Cells.Item(8, 12).Value2 = "Date de fin de l'effet indésirable"    # Short
$ws.Cells.Item(8, 13).Value2 = "Date de fin de l'effet indésirable"    # Definition
$ws.Cells.Item(8, 32).Value2 = "fr-lm-effet-indesirable.dateFin"       # Base Path

# Re-apply the same cell formatting/borders as the row above so the new
# row visually matches the rest of the table.
$ws.Range("A7:AJ7").Copy()
$ws.Range("A8:AJ8").PasteSpecial(-4122)  # xlPasteFormats
